$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to reflect the new Dataset-based naming scheme
$ws.Range("A1").Value = "PIDN_x"
$ws.Range("B1").Value = "DCDate_x"
$ws.Range("J1").Value = "_mp_merge"
$ws.Range("K1").Value = "_mp_diff_days"
$ws.Range("L1").Value = "_mp_abs_diff_days"

# Remove the now-unused _duplicates column entirely
$ws.Range("M1:M5").EntireColumn.Delete()

# Resize columns to (best) fit the new, longer/shorter header text
$ws.Range("A1").EntireColumn.ColumnWidth = 5.833333333333333
$ws.Range("J1").EntireColumn.ColumnWidth = 9.5
$ws.Range("K1").EntireColumn.ColumnWidth = 11.166666666666666
$ws.Range("L1").EntireColumn.ColumnWidth = 14.5
